$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entidade")
$ws.Activate()

# New shared string value in K2
$ws.Range("K2").Value = "CASA AMARELA"

# I2, J2 new values
$ws.Range("I2").Value = 6807000
$ws.Range("J2").Value = 10

# L2:O2 = "sim"
$ws.Range("L2").Value = "sim"
$ws.Range("M2").Value = "sim"
$ws.Range("N2").Value = "sim"
$ws.Range("O2").Value = "sim"

# Apply centered alignment to full A1:O2 range
$rng = $ws.Range("A1:O2")
$rng.HorizontalAlignment = -4108  # xlCenter
$rng.VerticalAlignment = -4108    # xlCenter

# Special font for I2
$ws.Range("I2").Font.Name = "Fira Code"
$ws.Range("I2").Font.Color = 3771391

# Selection
$ws.Range("E7").Select()
